$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.253.74"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").Value = "3.491.42"
$ws.Range("E3").Value = "  -0.30%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.51%  "

$ws.Range("D7").Value = "3.491.81"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("E10").Value = "  +2.01%  "

$ws.Range("E11").Value = "  -5.01%  "

$ws.Range("E12").Value = "  -2.72%  "

$ws.Range("D13").Value = "4.083.09"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.77%  "

$ws.Range("E15").Value = "  -5.31%  "

$ws.Range("D16").Value = "3.483.97"
$ws.Range("E16").Value = "  -0.62%  "

$ws.Range("D17").Value = "66.281.06"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.87%  "

$ws.Range("E20").Value = "  -3.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.594"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("D25").Value = "3.625.81"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E27").Value = "  -2.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "

$ws.Range("D35").Value = "3.469.22"
$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("E38").Value = "  -5.09%  "

$ws.Range("E39").Value = "  -3.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "169.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0859"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.880"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("E45").Value = "  -8.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.33%  "

$ws.Range("E48").Value = "  -9.87%  "

$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.942"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.71%  "
